# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: AD1 = "Wins", AE1 = "Losses", AF1 = "Ties" ---
# Copy the formatting (style) from the existing header cell A1 so the new
# header cells get the same bold/border/centered style (style index 1).
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows 2-50: AD = 92 (Wins), AE = 70 (Losses), AF = 0 (Ties) ---
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 92   # AD
    $ws.Cells.Item($r, 31).Value = 70   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
